$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.617.02"
$ws.Range("E2").Value = "  +3.59%  "
$ws.Range("D3").Value = "1.921.23"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.82"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.696"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.26"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.94"
$ws.Range("E9").Value = "  +10.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.367"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("E11").Value = "  +3.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0997"
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.46"
$ws.Range("E13").Value = "  +8.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.800"
$ws.Range("E14").Value = "  +5.59%  "
$ws.Range("D15").Value = "2.199.35"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.13"
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("D17").Value = "1.917.11"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "36.528.02"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.38"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("E20").Value = "  +4.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "251.12"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.29"
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("E23").Value = "  +5.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.92"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.80"
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.81"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("E31").Value = "  +6.51%  "
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("E33").Value = "  +4.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.35"
$ws.Range("E34").Value = "  +5.16%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0879"
$ws.Range("E36").Value = "  +25.99%  "
$ws.Range("E37").Value = "  -12.91%  "
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.46"
$ws.Range("E39").Value = "  +43.70%  "
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.49"
$ws.Range("E41").Value = "  +13.00%  "
$ws.Range("E42").Value = "  +5.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.24"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.12"
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("D45").Value = "1.348.43"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0811"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.46"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").Value = "2.097.69"
$ws.Range("E51").Value = "  +1.82%  "
